$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "A Multi-Phase Network Situational Awareness Cognitive Task Analysis"

$ws.Range("B3").Value = "Cyber scares and prophylactic policies: Crossnational evidence on the effect of cyberattacks on public support for surveillance"
$ws.Range("F3").Value = "Free Access"

$ws.Range("B4").Value = "The code not taken: China, the United States, and the future of cyber espionage"
$ws.Range("F4").Value = "Free Access"

$ws.Range("B5").Value = "On domains: Cyber and the practice of warfare"

$ws.Range("B6").Value = "Moving beyond the sanctuary paradigm: Canada must face up to the reality of a contested and dangerous space environment"

$ws.Range("B7").Value = "Responding to Uncertainty: The Importance of Covertness in Support for Retaliation to Cyber and Kinetic Attacks"
$ws.Range("F7").Value = "Open Access"

$ws.Range("B8").Value = "Cyclones in cyberspace: Information shaping and denial in the 2008 Russia–Georgia war"
$ws.Range("F8").Value = "Free Access"

$ws.Range("B9").Value = "Cyber and contentious politics: Evidence from the US radical environmental movement"
$ws.Range("F9").Value = "Free Access"

$ws.Range("B10").Value = "Digital Assays Part II: Digital Protein and Cell Assays"
$ws.Range("F10").Value = "Free Access"

$ws.Range("B11").Value = "Cyber-Flirting: Playing at Love on the Internet"
